# update database and change read_price algorithm
#
# The workbook tracks quarterly figures in columns E:N (10 quarters).
# This edit rolls the window forward by one quarter: the oldest quarter
# ("Q2 ended 1399/06") is dropped, every remaining quarter's column
# shifts one position to the left, and a new quarter ("Q4 ended 1401/12")
# is appended as the new column N - both for the quarter-header labels
# (rows 8 and 24) and for every data row beneath them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sequence of quarter labels for columns E..N (after the roll-forward)
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

# Update the two header rows (row 8 and row 24) that display the quarter labels
for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i   # column E = 5 .. column N = 14
    $ws.Cells.Item(8, $col).Value = $quarters[$i]
    $ws.Cells.Item(24, $col).Value = $quarters[$i]
}

# New values (already shifted one quarter to the left, with the new
# quarter's figure appended as column N) for every data row.
$rowValues = @{
    10 = @(-251, 1266, 442, 664, 3955, -1337, 1603, 2563, 3782, 2781)
    13 = @(0, 4568, 72, 642, 874, 2788, 5424, 1786, 1733, 730)
    14 = @(720, 1951, 1144, 2232, 2385, 2227, 2133, 2918, 604, 4220)
    15 = @(2967, 190, 892, 5088, 4581, 3188, 4384, 6526, 7155, 5092)
    16 = @(1034, 1741, 1762, 1827, 1515, 6217, 5693, 5835, 5848, 6398)
    17 = @(30216, 50602, 47897, 54168, 32774, 96197, 75108, 81557, 103872, 139144)
    19 = @(14991, 57531, 17384, 20540, 27943, 32211, 11063, 25112, 15560, 78760)
    20 = @(49677, 117849, 69593, 85161, 74027, 141491, 105408, 126297, 138554, 237125)
    26 = @(368, 173, 171, 180, 180, 165, 158, 158, 154, 152)
    27 = @(598, 796, 792, 790, 790, 805, 801, 810, 827, 829)
}

foreach ($row in $rowValues.Keys) {
    $vals = $rowValues[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i   # column E = 5 .. column N = 14
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

$wb.Save()
